# Apply the "automatic update of files" edit:
#  - Column C (Förändrad) for every data row (2-16) changes from 46077 to 46078.
#  - Rows 6-16 get their A (Beteckning), B (Datum) and G (Area ha) values
#    reshuffled to new positions (the underlying list was re-sorted/updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column C (Förändrad) for all data rows 2..16 ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46078
}

# --- 2. Rewrite rows 6..16 (columns A, B, G) with their new contents ---
$rowsData = @{
    6  = @("A 2593-2024",  45313.69204861111, 2.3)
    7  = @("A 12651-2022", 44641,             3.2)
    8  = @("A 8194-2025",  45708,             1.9)
    9  = @("A 5792-2024",  45335,             5.6)
    10 = @("A 13651-2023", 45006,             2.2)
    11 = @("A 50997-2025", 45946,             1.5)
    12 = @("A 35642-2023", 45147,             1.2)
    13 = @("A 7827-2026",  46062.63958333333, 2.1)
    14 = @("A 7814-2026",  46062.61388888889, 1.1)
    15 = @("A 28288-2023", 45099.6349537037,  0.5)
    16 = @("A 7333-2025",  45703.35899305555, 0.9)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
